$wb = $excel.ActiveWorkbook

# --- "inputs" sheet: rename parameter v_start -> starting_speed ---
$ws1 = $wb.Worksheets.Item("inputs")
$ws1.Range("A4").Value2 = "starting_speed"
$ws1.Range("C14").Select()

# --- "time_series" sheet: replace velocities column with a computed
#     running index in column A (times = 0, 1, 2, ... 20) and clear the
#     old "velocities" sample data in column B ---
$ws4 = $wb.Worksheets.Item("time_series")
$ws4.Range("B2:B22").ClearContents()
$ws4.Range("A3").Formula = "=A2 + 1"
$ws4.Range("A4:A22").Formula = "=A3 + 1"
$ws4.Range("K18").Select()
